$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the last row (RF_012 / Modificar_top_canciones) from the
#    "Requisitos funcionales del usuario administrador de la emisora"
#    table, and move the "_GoBack" bookmark (which lived in that row)
#    onto the now-empty paragraph that follows the table.
# ------------------------------------------------------------------

# Locate the requirements table whose last row holds RF_012 / Modificar_top_canciones.
$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    $lastRowText = $tbl.Rows.Item($tbl.Rows.Count).Range.Text
    if ($lastRowText -like "*Modificar_top_canciones*") {
        $targetTable = $tbl
        break
    }
}

if ($targetTable -ne $null) {
    # Drop the final row entirely.
    $targetTable.Rows.Item($targetTable.Rows.Count).Delete()

    # The paragraph right after the table is now a genuinely empty
    # paragraph (no runs at all). Adding a bookmark straight into a
    # zero-length range there is unreliable, so we briefly give the
    # paragraph a placeholder character, anchor the bookmark next to
    # it, and then remove the placeholder again.
    $afterPos = $targetTable.Range.End
    $placeholderRng = $d.Range($afterPos, $afterPos)
    $placeholderRng.InsertAfter("x")

    $bmRng = $d.Range($afterPos, $afterPos)
    $d.Bookmarks.Add("_GoBack", $bmRng)

    $placeholderRng2 = $d.Range($afterPos, $afterPos + 1)
    $placeholderRng2.Delete()
}

# ------------------------------------------------------------------
# 2) Drop the stale lastRenderedPageBreak cached on the final
#    "Anexos" heading run (re-assigning the paragraph text clears the
#    cached rendering marker while preserving the run formatting).
# ------------------------------------------------------------------

$paras = $d.Content.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "Anexos`r") {
        $p.Range.Text = "Anexos"
    }
}
